# Sync automático del tracker - actualiza resultados de predicciones completadas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sentAt = "2025-09-08 04:28:17"

# Row 74: Houston Dynamo vs Los Angeles Galaxy -> Draw (Fallo)
$ws.Range("L74").Value = "Completed"
$ws.Range("M74").Value = "Draw"
$ws.Range("N74").Value = "Fallo"
$ws.Range("O74").Value = -2
$ws.Range("P74").Value = -100
$ws.Range("Q74").Value = $sentAt

# Row 75: Chicago Fire vs New England Revolution -> Home Win (Acierto)
$ws.Range("L75").Value = "Completed"
$ws.Range("M75").Value = "Home Win"
$ws.Range("N75").Value = "Acierto"
$ws.Range("O75").Value = 1.65
$ws.Range("P75").Value = 75
$ws.Range("Q75").Value = $sentAt

# Row 76: St. Louis City vs FC Dallas -> Draw (Fallo)
$ws.Range("L76").Value = "Completed"
$ws.Range("M76").Value = "Draw"
$ws.Range("N76").Value = "Fallo"
$ws.Range("O76").Value = -1
$ws.Range("P76").Value = -100
$ws.Range("Q76").Value = $sentAt

# Row 77: CDS Tampico Madero vs Tapatío -> Home Win (Acierto)
$ws.Range("L77").Value = "Completed"
$ws.Range("M77").Value = "Home Win"
$ws.Range("N77").Value = "Acierto"
$ws.Range("O77").Value = 1.49
$ws.Range("P77").Value = 62
$ws.Range("Q77").Value = $sentAt

# Row 78: Tepatitlán vs Leones Negros UDG -> Draw (Fallo)
$ws.Range("L78").Value = "Completed"
$ws.Range("M78").Value = "Draw"
$ws.Range("N78").Value = "Fallo"
$ws.Range("O78").Value = -1.1
$ws.Range("P78").Value = -100
$ws.Range("Q78").Value = $sentAt
